# Updated cryptos list values to match the latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.755.82'
$ws.Range("E2").Value = '  +0.66%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''606.08'
$ws.Range("E5").Value = '  +1.93%  '
$ws.Range("D6").Value = '''141.44'
$ws.Range("E6").Value = '  +0.72%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.304.56'
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("E10").Value = '  +1.77%  '
$ws.Range("E11").Value = '  +3.75%  '
$ws.Range("E12").Value = '  +1.01%  '
$ws.Range("D13").Value = '''0.0000247'
$ws.Range("E13").Value = '  +0.26%  '
$ws.Range("D14").Value = '''34.90'
$ws.Range("E14").Value = '  +1.71%  '
$ws.Range("D15").Value = '3.852.12'
$ws.Range("E15").Value = '  +2.13%  '
$ws.Range("E16").Value = '  +0.81%  '
$ws.Range("D17").Value = '3.307.46'
$ws.Range("E17").Value = '  +1.92%  '
$ws.Range("D18").Value = '63.846.41'
$ws.Range("E18").Value = '  +0.79%  '
$ws.Range("E19").Value = '  +1.84%  '
$ws.Range("D20").Value = '''480.85'
$ws.Range("E20").Value = '  +1.51%  '
$ws.Range("D21").Value = '''14.02'
$ws.Range("E21").Value = '  -1.11%  '
$ws.Range("E22").Value = '  +1.02%  '
$ws.Range("D23").Value = '''8.00'
$ws.Range("E23").Value = '  +0.69%  '
$ws.Range("D24").Value = '''13.98'
$ws.Range("E24").Value = '  +6.22%  '
$ws.Range("E25").Value = '  +1.43%  '
$ws.Range("E27").Value = '  +1.73%  '
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("E29").Value = '  -0.32%  '
$ws.Range("E30").Value = '  +1.01%  '
$ws.Range("E31").Value = '  +1.25%  '
$ws.Range("D32").Value = '''28.75'
$ws.Range("E32").Value = '  +4.27%  '
$ws.Range("E33").Value = '  -1.14%  '
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("E35").Value = '  +1.36%  '
$ws.Range("D36").Value = '''6.06'
$ws.Range("E36").Value = '  +2.29%  '
$ws.Range("D37").Value = '''52.48'
$ws.Range("E37").Value = '  -0.42%  '
$ws.Range("D38").Value = '0.0₃0746'
$ws.Range("E38").Value = '  +5.17%  '
$ws.Range("E39").Value = '  +1.71%  '
$ws.Range("D40").Value = '3.112.79'
$ws.Range("E40").Value = '  +4.47%  '
$ws.Range("D41").Value = '''430.23'
$ws.Range("E41").Value = '  +1.80%  '
$ws.Range("E42").Value = '  +7.49%  '
$ws.Range("B43").Value = 'Cosmos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D43").Value = '''8.33'
$ws.Range("E43").Value = '  -0.64%  '
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").Value = '''2.74'
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").Value = '''0.265'
$ws.Range("E45").Value = '  -0.41%  '
$ws.Range("E46").Value = '  +2.75%  '
$ws.Range("D47").Value = '''36.80'
$ws.Range("E47").Value = '  +9.12%  '
$ws.Range("D48").Value = '''26.37'
$ws.Range("E48").Value = '  +1.88%  '
$ws.Range("D50").Value = '''126.40'
$ws.Range("E50").Value = '  +4.14%  '
$ws.Range("E51").Value = '  +0.08%  '
